# Apply cryptos list update (prices / 1h volume deltas) - Fri Aug 25 13:57:55 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.169.02"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.660.69"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.15"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5216"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2637"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06272"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.74"
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.465"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.661.96"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").Value = "1.890.93"
$ws.Range("E14").Value = "  -0.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5445"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").Value = "0.0₅8132"
$ws.Range("E16").Value = "  -1.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.90"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "26.193.91"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.592"
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.56"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.00"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.001"
$ws.Range("E23").Value = "  -4.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.71"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -2.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.247"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.16"
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.409"
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05937"
$ws.Range("E30").Value = "  -3.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.272"
$ws.Range("E31").Value = "  -1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.529"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.262"
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.575"
$ws.Range("E34").Value = "  -6.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9561"
$ws.Range("E35").Value = "  -4.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.421"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.770"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5655"
$ws.Range("E38").Value = "  -6.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01595"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.956"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8501"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.53"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "1.003.75"
$ws.Range("E44").Value = "  -7.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.58"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("E47").Value = "  -2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.013"
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4346"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05150"
$ws.Range("E51").Value = "  -0.99%  "
